$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 97.62036855931197
$ws.Range("H2").Value = 97.25341747731727
$ws.Range("I2").Value = 95.94290401183272

$ws.Range("G3").Value = 97.89434264495286
$ws.Range("H3").Value = 97.3195005698543
$ws.Range("I3").Value = 95.93138923106169

$ws.Range("G4").Value = 97.75204231889929
$ws.Range("H4").Value = 97.24207440014614
$ws.Range("I4").Value = 95.89315695024639

$ws.Range("G5").Value = 97.68724301390701
$ws.Range("H5").Value = 97.19735533611643
$ws.Range("I5").Value = 95.92185052869699

$ws.Range("G6").Value = 97.79240300414726
$ws.Range("H6").Value = 97.22416751310577
$ws.Range("I6").Value = 95.85427150582643
